# Update cryptos list (prices and 1h volume change) as of
# Wed Jan 10 03:23:46 UTC 2024.
#
# Most "Price" (column D) values are plain numeric-looking strings that
# must stay as TEXT (the source sheet stores everything as inline/shared
# strings). Excel's COM layer auto-coerces a plain numeric-looking string
# assigned via .Value into a real number, so for those cells we first force
# the cell's number format to Text ("@") before assigning the string. Values
# that contain extra separators (e.g. "46.213.86") are never parsed as a
# number by Excel, so they can be assigned directly without that treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $cellRef, $val) {
    $r = $sheet.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "46.213.86"
$ws.Range("E2").Value = "  -1.51%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.357.66"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4 - TetherUSD
Set-TextValue $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "301.85"
$ws.Range("E5").Value = "  +0.94%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "99.80"
$ws.Range("E6").Value = "  +1.44%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.571"
$ws.Range("E7").Value = "  -0.50%  "

# Row 8 - USDC
Set-TextValue $ws "D8" "1.00"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
Set-TextValue $ws "D9" "0.513"
$ws.Range("E9").Value = "  -3.10%  "

# Row 10 - Avalanche
Set-TextValue $ws "D10" "34.41"
$ws.Range("E10").Value = "  -3.82%  "

# Row 11 - Dogecoin
Set-TextValue $ws "D11" "0.0799"
$ws.Range("E11").Value = "  -0.02%  "

# Row 12 - Polkadot
Set-TextValue $ws "D12" "7.14"
$ws.Range("E12").Value = "  -2.67%  "

# Row 13 - TRON
Set-TextValue $ws "D13" "0.104"
$ws.Range("E13").Value = "  -0.30%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.717.84"
$ws.Range("E14").Value = "  +2.05%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "2.358.44"
$ws.Range("E15").Value = "  +2.18%  "

# Rows 16 & 17 swap places: Chainlink/Polygon order flips in the ranking
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws "D16" "0.810"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D17" "13.65"
$ws.Range("E17").Value = "  -2.32%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "46.122.30"
$ws.Range("E18").Value = "  -1.44%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue $ws "D19" "12.76"
$ws.Range("E19").Value = "  -2.83%  "

# Row 20 - ShibaInu (price unchanged, only volume changes)
$ws.Range("E20").Value = "  +2.65%  "

# Row 21 - Uniswap
Set-TextValue $ws "D21" "6.07"
$ws.Range("E21").Value = "  -0.97%  "

# Row 22 - Litecoin
Set-TextValue $ws "D22" "67.36"
$ws.Range("E22").Value = "  +0.72%  "

# Row 23 - BitcoinCash
Set-TextValue $ws "D23" "246.97"
$ws.Range("E23").Value = "  -0.61%  "

# Row 24 - PancakeSwap
Set-TextValue $ws "D24" "2.85"
$ws.Range("E24").Value = "  -2.68%  "

# Row 25 - Dai (price unchanged, only volume changes)
$ws.Range("E25").Value = "  -0.01%  "

# Row 26 - ImmutableX (price unchanged, only volume changes)
$ws.Range("E26").Value = "  -3.14%  "

# Row 27 - InjectiveProtocol
Set-TextValue $ws "D27" "39.68"
$ws.Range("E27").Value = "  -7.72%  "

# Row 28 - Toncoin (price unchanged, only volume changes)
$ws.Range("E28").Value = "  -2.69%  "

# Row 29 - Cosmos
Set-TextValue $ws "D29" "9.79"
$ws.Range("E29").Value = "  -0.49%  "

# Row 30 - EthereumClassic
Set-TextValue $ws "D30" "21.02"
$ws.Range("E30").Value = "  +4.19%  "

# Row 31 - LidoDAOToken
Set-TextValue $ws "D31" "3.72"
$ws.Range("E31").Value = "  +19.63%  "

# Row 32 - WEMIXToken
Set-TextValue $ws "D32" "2.78"
$ws.Range("E32").Value = "  +5.76%  "

# Row 33 - Filecoin
Set-TextValue $ws "D33" "5.55"
$ws.Range("E33").Value = "  -3.39%  "

# Row 34 - Monero
Set-TextValue $ws "D34" "146.08"
$ws.Range("E34").Value = "  -0.95%  "

# Row 35 - Hedera
Set-TextValue $ws "D35" "0.0775"
$ws.Range("E35").Value = "  -2.84%  "

# Row 36 - Kaspa
Set-TextValue $ws "D36" "0.112"
$ws.Range("E36").Value = "  -0.38%  "

# Row 37 - ARBITRUM
Set-TextValue $ws "D37" "1.88"
$ws.Range("E37").Value = "  +4.78%  "

# Row 38 - Stellar
Set-TextValue $ws "D38" "0.117"
$ws.Range("E38").Value = "  -1.99%  "

# Row 39 - Celestia
Set-TextValue $ws "D39" "14.97"
$ws.Range("E39").Value = "  -5.51%  "

# Row 40 - RenderToken
Set-TextValue $ws "D40" "3.95"
$ws.Range("E40").Value = "  -0.80%  "

# Row 41 - VeChain
Set-TextValue $ws "D41" "0.0301"
$ws.Range("E41").Value = "  -2.32%  "

# Row 42 - NEARProtocol
Set-TextValue $ws "D42" "3.23"
$ws.Range("E42").Value = "  -6.34%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.883.30"
$ws.Range("E43").Value = "  +2.47%  "

# Row 44 - FirstDigitalUSD
Set-TextValue $ws "D44" "0.999"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45 - BitcoinSV
Set-TextValue $ws "D45" "93.42"
$ws.Range("E45").Value = "  +2.45%  "

# Row 46 - Stacks
Set-TextValue $ws "D46" "1.80"
$ws.Range("E46").Value = "  -9.76%  "

# Row 47 - Algorand (price unchanged, only volume changes)
$ws.Range("E47").Value = "  -6.08%  "

# Row 48 - FraxShare
Set-TextValue $ws "D48" "8.26"
$ws.Range("E48").Value = "  +3.47%  "

# Row 49 - Aave
Set-TextValue $ws "D49" "97.84"
$ws.Range("E49").Value = "  +0.53%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.588.78"
$ws.Range("E50").Value = "  +1.93%  "

# Row 51 - ordi
Set-TextValue $ws "D51" "68.97"
$ws.Range("E51").Value = "  -9.42%  "
